$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("Sheet1" - the main Log sheet) ----
$ws1 = $wb.Worksheets.Item(1)

# Row 12: "Check" -> "check"
$ws1.Range("G12").Value = "check"

# Row 13: add "check" in G13, update H13 text
$ws1.Range("G13").Value = "check"
$ws1.Range("H13").Value = "July 26 - July 28"

# Row 14: add Array 101 note + next date range
$ws1.Range("F14").Value = "Array 101"
$ws1.Range("H14").Value = "July 29 - July "

# Row 21: new log entry for 44040 (July 28 2020)
$ws1.Range("A21").Value = 44040
$ws1.Range("C21").Value = "P394, P542, P733, P841, P43"
$ws1.Range("D21").Value = 4
$ws1.Range("E21").Value = "Cook"

# ---- Sheet2 ("一刷看答案") ----
$ws2 = $wb.Worksheets.Item(2)

# A2/A3 need the same date style as A1 (s="9"); copy formats down first
$ws2.Range("A1").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

$ws2.Range("A2").Value = 44040
$ws2.Range("B2").Value = 394
$ws2.Range("C2").Value = 733
$ws2.Range("D2").Value = 841
$ws2.Range("E2").Value = 43

$ws2.Range("A3").Value = 44041

# update selection on sheet2 first (becomes the active sheet momentarily)
$null = $ws2.Range("B3").Select()

# finally, re-select on Sheet1 so it ends up the active/visible tab again
$null = $ws1.Range("C21").Select()
